$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102:212 down to 103:213.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new record's data.
$ws.Cells.Item(102, 1).Value = 3
$ws.Cells.Item(102, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(102, 3).Value = "Coquimbo"
$ws.Cells.Item(102, 4).Value = 44494
$ws.Cells.Item(102, 5).Value = 5
$ws.Cells.Item(102, 6).Value = 100114013
$ws.Cells.Item(102, 7).Value = "Zanahoria"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 190
$ws.Cells.Item(102, 11).Value = 8000
$ws.Cells.Item(102, 12).Value = 8000
$ws.Cells.Item(102, 13).Value = 8000
$ws.Cells.Item(102, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(102, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(102, 16).Value = 400
$ws.Cells.Item(102, 17).Value = 20
$ws.Cells.Item(102, 18).Value = "Hortaliza"
